# ProcessBook.docx edit script
# 1) Adds a new sub-bullet paragraph (ilvl=1) right after the "Design Evolution"
#    bullet, describing the transition-design decision, matching the style of
#    the existing "Initially, we made a single scatterplot..." sub-bullet.
# 2) Moves the hidden "_GoBack" bookmark from the end of the "Design Evolution"
#    bullet to the end of the "Data:" bullet (reflecting where the edit
#    session's cursor ended up).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Insert the new paragraph after "Design Evolution" bullet (para 10)
# ---------------------------------------------------------------------------

# Locate the two anchor paragraphs by their known text.
$designEvolutionPara = $null
$analysisPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("Design Evolution:")) { $designEvolutionPara = $p }
    if ($t.StartsWith("Analysis:")) { $analysisPara = $p }
}

# Find the existing "Initially, we made a single scatterplot..." sub-bullet
# paragraph - it already has the exact pPr/rPr formatting (ilvl=1 numbered
# sub-bullet, Times New Roman, color 111111) that the new paragraph needs, so
# copy it wholesale and then swap in the new text.
$templatePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Initially, we made a single scatterplot")) {
        $templatePara = $p
    }
}

$templatePara.Range.Copy()

$insertPos = $designEvolutionPara.Range.End
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.Paste()

# The pasted paragraph is now the one right after Design Evolution and before
# Analysis. Re-find it by position (paragraph immediately before Analysis).
$newPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Initially, we made a single scatterplot") -and $p.Range.Start -ne $templatePara.Range.Start) {
        $newPara = $p
    }
}

# Replace its text (minus trailing paragraph mark) with the new sentence,
# collapsing it down to a single run just like the target content.
$newText = "When we were looking at how to transition data when we were changing positions, at first, the interaction had all old points rising toward the top left point while fading out, whereas the new points were falling from the top left point while fading in. After discussion, we decided that it would be a less confusing transition to just have old points fade out and new points fade in."
$bodyRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$bodyRange.Text = $newText

# ---------------------------------------------------------------------------
# Step 2: Move the "_GoBack" bookmark to the end of the "Data:" bullet
# ---------------------------------------------------------------------------

# Re-find the Data: paragraph (indices shifted after the insert above).
$dataPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Data:")) { $dataPara = $p }
}

# Remove the old bookmark (currently at the end of the Design Evolution bullet).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Adding a bookmark collapsed exactly on a paragraph mark is unreliable, so
# insert a temporary placeholder character immediately before the paragraph
# mark, anchor the bookmark just before that placeholder (a safe, mid-run
# position), and then remove the placeholder - leaving a zero-width bookmark
# sitting right at the end of the paragraph's text, before the mark.
$pilcrowPos = $dataPara.Range.End - 1
$placeholderRange = $d.Range($pilcrowPos, $pilcrowPos)
$placeholderRange.InsertBefore("X")

$bookmarkRange = $d.Range($pilcrowPos, $pilcrowPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$placeholderCharRange = $d.Range($pilcrowPos, $pilcrowPos + 1)
$placeholderCharRange.Delete()
